$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("candidate_preferences")

# The "Jerry" entry in A4 is actually "Ginger" -- fix the candidate name.
# This both rewrites the cell's value and appends "Ginger" as a new shared
# string (sharedStrings uniqueCount 8 -> 9).
$ws.Range("A4").Value = "Ginger"

# Move/leave the active selection on the corrected cell.
$ws.Range("A4").Select() | Out-Null
